# Update PLC data 2025-10-13 14:16:38
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7208
$ws.Range("C3").Value = 179969
$ws.Range("C4").Value = 169926
$ws.Range("C8").Value = 65.11
